# lesson3.docx edit: the single run holding the GitHub URL gets split
# into three runs (with the existing "_GoBack" bookmark re-anchored in
# between runs 2 and 3) and the "l" of "lesson3" becomes a capital "L".
#
#   before: "https://github.com/RefrigeratorUA/GoITLib/tree/master/lesson3"
#   after : "https://github.com/Refr" | "igeratorUA/GoITLib/tree/master/L" | <bookmark> | "esson3"

$d = $word.ActiveDocument

$oldUrl = "https://github.com/RefrigeratorUA/GoITLib/tree/master/lesson3"
$run1   = "https://github.com/Refr"
$run2   = "igeratorUA/GoITLib/tree/master/L"
$run3   = "esson3"

# Locate the run holding the URL and capture where it starts.
$hit = $d.Content
$hit.Find.Execute($oldUrl, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$base = $hit.Start

$splitAt1 = $base + $run1.Length            # boundary between run 1 / run 2
$splitAt2 = $base + $run1.Length + $run2.Length   # boundary between run 2 / run 3 (bookmark goes here)

# 1) Fix the casing: "l" -> "L" right before "esson3".
$d.Range($splitAt2 - 1, $splitAt2).Text = "L"

# 2) Force a run break between "...Refr" and "igerator..." by dropping a
#    throwaway bookmark at the boundary and immediately deleting it again
#    (adding/removing a bookmark splits the underlying run but leaves no
#    formatting residue behind).
$d.Bookmarks.Add("zzTmpSplit", $d.Range($splitAt1, $splitAt1)) | Out-Null
$d.Bookmarks("zzTmpSplit").Delete()

# 3) Re-anchor the existing "_GoBack" bookmark from the end of the text to
#    the boundary between run 2 ("...master/L") and run 3 ("esson3"); this
#    both moves the bookmark and produces the final run split.
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($splitAt2, $splitAt2)) | Out-Null

Write-Host "Final paragraph text:" $d.Paragraphs(1).Range.Text
